# The sheet holds monthly index values in columns A (month label, "YYYY-MM")
# and D (numeric value), one row per month, grouped into contiguous 12-row
# blocks per year starting at row 2 (rows 2-13 = 2014, 14-25 = 2015,
# 26-37 = 2016, 38-49 = 2017).
#
# The edit re-orders the months *within each year block* so that
# Oct, Nov, Dec come first, followed by Jan..Sep (a left-rotate-by-9 /
# right-rotate-by-3 of the 12 rows). Column B/C are untouched (they are
# blank for every data row already). We read the existing values straight
# off the sheet so the script is driven by whatever data is actually
# present, then write the rotated order back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$rowsPerYear = 12
$numYears = 4

for ($yearIndex = 0; $yearIndex -lt $numYears; $yearIndex++) {
    $blockStart = $firstDataRow + ($yearIndex * $rowsPerYear)

    # Snapshot the 12 (month, value) pairs for this year before writing
    # anything back, so overlapping reads/writes inside the same block
    # can't clobber each other.
    $months = @()
    $values = @()
    for ($i = 0; $i -lt $rowsPerYear; $i++) {
        $r = $blockStart + $i
        # NOTE: call .Value() (not the bare property) to force the getter.
        $months += , ($ws.Cells.Item($r, 1).Value())
        $values += , ($ws.Cells.Item($r, 4).Value())
    }

    # Rotate: last 3 (Oct, Nov, Dec => indices 9,10,11) move to the front,
    # followed by the first 9 (Jan..Sep => indices 0..8).
    $order = @(9, 10, 11, 0, 1, 2, 3, 4, 5, 6, 7, 8)

    for ($i = 0; $i -lt $rowsPerYear; $i++) {
        $r = $blockStart + $i
        $srcIndex = $order[$i]
        $ws.Cells.Item($r, 1).Value = $months[$srcIndex]
        $ws.Cells.Item($r, 4).Value = $values[$srcIndex]
    }
}
